# Fleet10and11_MRFSSDockside_IndexData.xlsx -- "All tables for indices in"
#
# Summary of the edit:
#  - NorthAIC / SouthAIC: AIC numbers >=1000 are re-entered as comma-formatted
#    text (e.g. 1481 -> "1,481"), and "Area_X" is re-escaped to "Area\_X".
#  - BothFilter: rebuilt to show BOTH the north and south filtering cascades
#    (previously only NorthFilter/SouthFilter held this data; BothFilter had
#    a half-finished/placeholder table). Big numbers are comma-formatted text,
#    the "drifts" counts stay numeric. "NA" rows separate the two cascades.
#  - BothFilter becomes the active/selected sheet (was SouthIndex before).
#  - Minor selection-cell housekeeping on the other sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# NorthAIC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("NorthAIC")

$ws.Range("A4").Value = "Year + Region + Area\_X"
$ws.Range("A5").Value = "Year + Region + Area\_X + Wave"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1,481"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1,429"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "1,403"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "1,397"

$ws.Columns.Item(1).ColumnWidth = 27.42
$ws.Range("A5").Select()

# ---------------------------------------------------------------------------
# SouthAIC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SouthAIC")

$ws.Range("A4").Value = "Year + Wave + Area\_X"
$ws.Range("A5").Value = "Year + Wave + Area\_X + SubRegion"

$ws.Columns.Item(1).ColumnWidth = 52.25
$ws.Range("A5").Select()

# ---------------------------------------------------------------------------
# BothFilter -- full rewrite
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BothFilter")

# Header row (Filter / Drifts / PositiveDrifts) is unchanged.

# Rows 2-4: "All data" cascade, now shared with totals as comma text.
$ws.Range("B2:C4").NumberFormat = "@"
$ws.Range("A2").Value = "All data"
$ws.Range("B2").Value = "10,392"
$ws.Range("C2").Value = "1,061"
$ws.Range("A3").Value = "Remove north of Cape Mendocino"
$ws.Range("B3").Value = "10,327"
$ws.Range("C3").Value = "1,061"
$ws.Range("A4").Value = "Remove trips targetting offshore species"
$ws.Range("B4").Value = "10,122"
$ws.Range("C4").Value = "1,061"

# Row 5: separator
$ws.Range("B5:C5").NumberFormat = "@"
$ws.Range("A5").Value = "NA"
$ws.Range("B5").Value = "NA"
$ws.Range("C5").Value = "NA"

# Rows 6-9: northern filtering cascade (mirrors NorthFilter sheet)
$ws.Range("A6").Value = "Start northern filtering"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "2,788"
$ws.Range("C6").Value = 620
$ws.Range("C6").NumberFormat = "@"

$ws.Range("A7").Value = "Remove species that never co-occurand  not present in at least 1% of all"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "2,788"
$ws.Range("C7").Value = 620
$ws.Range("C7").NumberFormat = "@"

$ws.Range("A8").Value = "Stephens-MacCall filter (keep all positives)"
$ws.Range("B8").Value = 806
$ws.Range("B8").NumberFormat = "@"
$ws.Range("C8").Value = 620
$ws.Range("C8").NumberFormat = "@"

$ws.Range("A9").Value = "Stephens-MacCall filter (keep only above threshold)"
$ws.Range("B9").Value = 623
$ws.Range("B9").NumberFormat = "@"
$ws.Range("C9").Value = 437
$ws.Range("C9").NumberFormat = "@"

# Row 10: separator
$ws.Range("B10:C10").NumberFormat = "@"
$ws.Range("A10").Value = "NA"
$ws.Range("B10").Value = "NA"
$ws.Range("C10").Value = "NA"

# Rows 11-14: southern filtering cascade (mirrors SouthFilter sheet)
$ws.Range("A11").Value = "Start southern filtering"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "7,334"
$ws.Range("C11").Value = 441
$ws.Range("C11").NumberFormat = "@"

$ws.Range("A12").Value = "Remove species that never co-occurand  not present in at least 1% of all"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "7,334"
$ws.Range("C12").Value = 441
$ws.Range("C12").NumberFormat = "@"

$ws.Range("A13").Value = "Stephens-MacCall filter (keep all positives)"
$ws.Range("B13").Value = 687
$ws.Range("B13").NumberFormat = "@"
$ws.Range("C13").Value = 441
$ws.Range("C13").NumberFormat = "@"

$ws.Range("A14").Value = "Stephens-MacCall filter (keep only above threshold)"
$ws.Range("B14").Value = 430
$ws.Range("B14").NumberFormat = "@"
$ws.Range("C14").Value = 184
$ws.Range("C14").NumberFormat = "@"

# ---------------------------------------------------------------------------
# NorthFilter -- selection housekeeping only
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("NorthFilter")
$ws.Range("A2:C5").Select()

# ---------------------------------------------------------------------------
# SouthFilter -- selection housekeeping only
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SouthFilter")
$ws.Range("B3").Select()

# ---------------------------------------------------------------------------
# BothFilter becomes the active sheet / tab (was SouthIndex before) -- must
# be the LAST Select() so it "wins" the workbook's activeTab/tabSelected.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BothFilter")
$ws.Range("C5").Select()
